# Commit: "Fruta / hortaliza, semanal"
# A new weekly price-report row is inserted at row 54 (pushing the
# existing rows 54-117 down to 55-118); the new row holds the latest
# observation for this market/product.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 54 - this shifts rows 54:117 down to 55:118 and
# copies the formatting (incl. the date number format on column D) from
# the row above, matching Excel's default insert behaviour.
$ws.Rows.Item(54).Insert()

$ws.Range("A54").Value = 8
$ws.Range("B54").Value = "Terminal La Palmera de La Serena"
$ws.Range("C54").Value = "Coquimbo"
$ws.Range("D54").Value = "2022-01-28"
$ws.Range("E54").Value = 4
$ws.Range("F54").Value = 100112001
$ws.Range("G54").Value = "Berenjena"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 500
$ws.Range("K54").Value = 8000
$ws.Range("L54").Value = 9000
$ws.Range("M54").Value = 8500
$ws.Range("N54").Value = "$/caja 50 unidades"
$ws.Range("O54").Value = "Región de Arica y Parinacota"
$ws.Range("P54").Value = 170
$ws.Range("Q54").Value = 50
$ws.Range("R54").Value = "Hortaliza"
